$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.733.43'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.546.69'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.24'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.92'
$ws.Range("E6").Value = '  +5.48%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.32'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.938.96'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.86'
$ws.Range("E15").Value = '  +5.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.523.57'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.839'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.762.85'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.40'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0956'
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.19'
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.22'
$ws.Range("E23").Value = '  -2.25%  '
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.07'
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.55'
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.05'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.36'
$ws.Range("E29").Value = '  -2.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.15'
$ws.Range("E30").Value = '  -2.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.66'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.75'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0805'
$ws.Range("E33").Value = '  +2.16%  '
$ws.Range("E34").Value = '  -2.69%  '
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("E36").Value = '  -2.80%  '
$ws.Range("E37").Value = '  +6.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.29'
$ws.Range("E38").Value = '  -6.03%  '
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("E40").Value = '  -0.51%  '
$ws.Range("E41").Value = '  +10.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.45'
$ws.Range("E42").Value = '  +3.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.29'
$ws.Range("E44").Value = '  +1.10%  '
$ws.Range("E45").Value = '  -1.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.981.00'
$ws.Range("E46").Value = '  -0.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.03'
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.793.20'
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.73'
$ws.Range("E51").Value = '  -1.41%  '
